# Applies the text edits on slide 2 (table of contents) and slide 18
# (results discussion) described by the commit diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 18: paragraph 3 of the body placeholder (shape 1)
#   "그러나, 새로운 손글씨 인식에서 다소 아쉬운 성능을 보입니다."
# becomes
#   "그러나, 새로운 글씨체의 손글씨 인식에서 다소 아쉬운 성능을 보입니다."
# with "손글씨" keeping its own (misspelling-flagged) run, a brand-new
# "글씨체의 " run ahead of it, and the trailing sentence split so a lone
# space becomes its own run before "인식에서 ...".
# ---------------------------------------------------------------------------
$slide18 = $p.Slides.Item(18)
$shape18 = $slide18.Shapes.Item(1)
$para18 = $shape18.TextFrame.TextRange.Paragraphs(3, 1)

# Step 1: grow the trailing space after "새로운" into " 글씨체의 " - this
# temporarily eats the space that belongs to "새로운 ".
$ins1 = $para18.Characters(9, 1)
$ins1.Text = " 글씨체의 "

# Step 2: re-carve "새로운 " (now spanning into the newly inserted text)
# back into its own run so the boundary sits where the diff expects it.
$fix1 = $para18.Characters(6, 4)
$fix1.Text = "새로운 "

# Step 3: split the long trailing run so the leading space in
# " 인식에서 다소 아쉬운 성능을 보입니다" becomes its own run.
$fix2 = $para18.Characters(18, 1)
$fix2.Text = " "

# ---------------------------------------------------------------------------
# Slide 2: table-of-contents placeholder (shape 2)
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange

# Paragraph 3: "머신 러닝" -> "딥" / "러닝"
$para3 = $tr2.Paragraphs(3, 1)
$word3 = $para3.Characters(1, 3)
$word3.Text = "딥"

# Paragraph 4: "생존율 예측 프로그램" -> "숫자 " / "손글씨" / " 인식" / " " / "프로그램"
$para4 = $tr2.Paragraphs(4, 1)
$whole4 = $para4.Characters(1, 11)
$whole4.Text = "숫자 손글씨 인식 프로그램"

$b1 = $para4.Characters(1, 3)
$b1.Text = $b1.Text
$b2 = $para4.Characters(4, 3)
$b2.Text = $b2.Text
$b3 = $para4.Characters(7, 3)
$b3.Text = $b3.Text
$b4 = $para4.Characters(10, 1)
$b4.Text = $b4.Text
